# Efna3-Epha2.xlsx update: refresh TPM-derived NATMI metrics and add the new
# MuSCs sending-cluster rows (Efna3 -> Epha2) for YoungD0.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna3"
$ws.Range("C2").Value = "Epha2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.046374
$ws.Range("H2").Value = 0.139122
$ws.Range("I2").Value = 0.6592866045237633
$ws.Range("J2").Value = 0.6592866045237632
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.757543333333333
$ws.Range("N2").Value = 23.27263
$ws.Range("O2").Value = 0.4040769763164727
$ws.Range("P2").Value = 0.4040769763164727
$ws.Range("Q2").Value = 0.35974831454
$ws.Range("R2").Value = 3.23773483086
$ws.Range("S2").Value = 0.2664025376819164
$ws.Range("T2").Value = 0.2664025376819164

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna3"
$ws.Range("C3").Value = "Epha2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.046374
$ws.Range("H3").Value = 0.139122
$ws.Range("I3").Value = 0.6592866045237633
$ws.Range("J3").Value = 0.6592866045237632
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.312365
$ws.Range("N3").Value = 0.937095
$ws.Range("O3").Value = 0.01627055103446774
$ws.Range("P3").Value = 0.01627055103446774
$ws.Range("Q3").Value = 0.01448561451
$ws.Range("R3").Value = 0.13037053059
$ws.Range("S3").Value = 0.01072695634524484
$ws.Range("T3").Value = 0.01072695634524484

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efna3"
$ws.Range("C4").Value = "Epha2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.046374
$ws.Range("H4").Value = 0.139122
$ws.Range("I4").Value = 0.6592866045237633
$ws.Range("J4").Value = 0.6592866045237632
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 11.12827366666667
$ws.Range("N4").Value = 33.384821
$ws.Range("O4").Value = 0.5796524726490594
$ws.Range("P4").Value = 0.5796524726490595
$ws.Range("Q4").Value = 0.516062563018
$ws.Range("R4").Value = 4.644563067162
$ws.Range("S4").Value = 0.3821571104966019
$ws.Range("T4").Value = 0.3821571104966019

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Efna3"
$ws.Range("C5").Value = "Epha2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.02396566666666667
$ws.Range("H5").Value = 0.071897
$ws.Range("I5").Value = 0.3407133954762367
$ws.Range("J5").Value = 0.3407133954762367
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 7.757543333333333
$ws.Range("N5").Value = 23.27263
$ws.Range("O5").Value = 0.4040769763164727
$ws.Range("P5").Value = 0.4040769763164727
$ws.Range("Q5").Value = 0.1859146976788889
$ws.Range("R5").Value = 1.67323227911
$ws.Range("S5").Value = 0.1376744386345563
$ws.Range("T5").Value = 0.1376744386345563

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Efna3"
$ws.Range("C6").Value = "Epha2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.02396566666666667
$ws.Range("H6").Value = 0.071897
$ws.Range("I6").Value = 0.3407133954762367
$ws.Range("J6").Value = 0.3407133954762367
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.312365
$ws.Range("N6").Value = 0.937095
$ws.Range("O6").Value = 0.01627055103446774
$ws.Range("P6").Value = 0.01627055103446774
$ws.Range("Q6").Value = 0.007486035468333333
$ws.Range("R6").Value = 0.06737431921500001
$ws.Range("S6").Value = 0.005543594689222901
$ws.Range("T6").Value = 0.005543594689222901

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Efna3"
$ws.Range("C7").Value = "Epha2"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.02396566666666667
$ws.Range("H7").Value = 0.071897
$ws.Range("I7").Value = 0.3407133954762367
$ws.Range("J7").Value = 0.3407133954762367
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 11.12827366666667
$ws.Range("N7").Value = 33.384821
$ws.Range("O7").Value = 0.5796524726490594
$ws.Range("P7").Value = 0.5796524726490595
$ws.Range("Q7").Value = 0.2666964972707778
$ws.Range("R7").Value = 2.400268475437
$ws.Range("S7").Value = 0.1974953621524575
$ws.Range("T7").Value = 0.1974953621524575
